$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Copy-CellTo([string]$srcRef, [string]$dstRef) {
    # Bring formatting in line with the source cell first ...
    $ws.Range($srcRef).Copy() | Out-Null
    $ws.Range($dstRef).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    # ... then copy the value itself (preserves text type, e.g. avoids a
    # "01/01/2012" string being reinterpreted as a date serial).
    $ws.Range($srcRef).Copy() | Out-Null
    $ws.Range($dstRef).PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# --- Relabel / move existing text into its new cells (dependency-safe order) ---
Copy-CellTo "B13" "B10"
Copy-CellTo "C13" "C10"
Copy-CellTo "A14" "A13"
Copy-CellTo "A15" "A14"
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()
Copy-CellTo "A16" "A15"
Copy-CellTo "B8" "B15"
Copy-CellTo "C8" "C15"
Copy-CellTo "A17" "A16"
$ws.Range("B16").Clear()
$ws.Range("C16").Clear()
Copy-CellTo "A18" "A17"
Copy-CellTo "A19" "A18"
Copy-CellTo "B13" "B18"
$ws.Range("B13").Value = "Semestral"
Copy-CellTo "C13" "C18"
$ws.Range("C13").Value = "Semestral"
Copy-CellTo "A20" "A19"
Copy-CellTo "A21" "A20"
Copy-CellTo "A22" "A21"
Copy-CellTo "A23" "A22"
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Range("A23").Clear()
Copy-CellTo "B24" "B23"
Copy-CellTo "C24" "C23"
Copy-CellTo "B25" "B24"
Copy-CellTo "C25" "C24"
$ws.Range("B25").Clear()
$ws.Range("C25").Clear()

# --- Remove the now-unused last row entirely ---
$ws.Rows("25:25").Delete()

# --- Row heights to match the new content layout ---
$ws.Rows("13:13").RowHeight = 60
$ws.Rows("15:15").RowHeight = 120
$ws.Rows("17:17").RowHeight = 15
$ws.Rows("18:18").RowHeight = 60
$ws.Rows("21:21").RowHeight = 120
$ws.Rows("22:22").RowHeight = 15
$ws.Rows("23:23").RowHeight = 30
